$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "79.769.06"
$ws.Range("E2").Value = "  +4.83%  "
$ws.Range("D3").Value = "3.215.40"
$ws.Range("E3").Value = "  +6.19%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "'210.87"
$ws.Range("E5").Value = "  +7.19%  "
$ws.Range("D6").Value = "'639.49"
$ws.Range("E6").Value = "  +3.18%  "
$ws.Range("D7").Value = "'0.265"
$ws.Range("E7").Value = "  +29.83%  "
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("D9").Value = "'0.601"
$ws.Range("E9").Value = "  +9.68%  "
$ws.Range("D10").Value = "3.214.06"
$ws.Range("E10").Value = "  +6.19%  "
$ws.Range("D11").Value = "'0.601"
$ws.Range("E11").Value = "  +36.97%  "
$ws.Range("D12").Value = "'0.0000271"
$ws.Range("E12").Value = "  +42.05%  "
$ws.Range("E13").Value = "  +3.53%  "
$ws.Range("D14").Value = "'5.44"
$ws.Range("E14").Value = "  +4.06%  "
$ws.Range("D15").Value = "3.807.59"
$ws.Range("E15").Value = "  +6.09%  "
$ws.Range("D16").Value = "'32.71"
$ws.Range("E16").Value = "  +13.57%  "
$ws.Range("D17").Value = "79.646.43"
$ws.Range("E17").Value = "  +4.69%  "
$ws.Range("D18").Value = "3.203.66"
$ws.Range("E18").Value = "  +5.74%  "
$ws.Range("D19").Value = "'14.67"
$ws.Range("E19").Value = "  +9.22%  "
$ws.Range("D20").Value = "'3.04"
$ws.Range("E20").Value = "  +29.20%  "
$ws.Range("D21").Value = "'9.44"
$ws.Range("E21").Value = "  +5.66%  "
$ws.Range("D22").Value = "'447.25"
$ws.Range("E22").Value = "  +17.51%  "
$ws.Range("D23").Value = "'5.28"
$ws.Range("E23").Value = "  +20.71%  "
$ws.Range("E24").Value = "  +12.59%  "
$ws.Range("E25").Value = "  +6.30%  "
$ws.Range("D26").Value = "'77.90"
$ws.Range("E26").Value = "  +7.43%  "
$ws.Range("E27").Value = "  +12.59%  "
$ws.Range("E28").Value = "  +0.17%  "
$ws.Range("D29").Value = "'0.0000127"
$ws.Range("E29").Value = "  +18.61%  "
$ws.Range("D30").Value = "'9.25"
$ws.Range("E30").Value = "  +12.29%  "
$ws.Range("D31").Value = "'0.997"
$ws.Range("E31").Value = "  -0.15%  "
$ws.Range("D32").Value = "'564.95"
$ws.Range("E32").Value = "  +14.62%  "
$ws.Range("E33").Value = "  +10.21%  "
$ws.Range("E34").Value = "  +32.14%  "
$ws.Range("E35").Value = "  +6.53%  "
$ws.Range("E36").Value = "  +13.19%  "
$ws.Range("E37").Value = "  +19.36%  "
$ws.Range("D38").Value = "'1.00"
$ws.Range("E38").Value = "  +0.04%  "
$ws.Range("E39").Value = "  +9.80%  "
$ws.Range("D40").Value = "'163.25"
$ws.Range("E40").Value = "  +0.78%  "
$ws.Range("B41").Value = "RenderToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D41").Value = "'5.75"
$ws.Range("E41").Value = "  +12.75%  "
$ws.Range("B42").Value = "WhiteBITCoin"
$ws.Range("C42").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D42").Value = "'20.28"
$ws.Range("E42").Value = "  +1.32%  "
$ws.Range("D43").Value = "'195.01"
$ws.Range("E43").Value = "  +2.12%  "
$ws.Range("E44").Value = "  +0.01%  "
$ws.Range("D45").Value = "'1.83"
$ws.Range("E45").Value = "  +12.24%  "
$ws.Range("D46").Value = "'2.72"
$ws.Range("E46").Value = "  +13.07%  "
$ws.Range("D47").Value = "'0.804"
$ws.Range("E47").Value = "  +4.49%  "
$ws.Range("D48").Value = "'1.34"
$ws.Range("E48").Value = "  +8.29%  "
$ws.Range("D49").Value = "'43.03"
$ws.Range("E49").Value = "  +4.10%  "
$ws.Range("E50").Value = "  +11.71%  "
$ws.Range("D51").Value = "'25.97"
$ws.Range("E51").Value = "  +17.04%  "
